$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 15.66920033333333
$ws.Range("H2").Value = 47.00760099999999
$ws.Range("I2").Value = 0.2925937299273087
$ws.Range("J2").Value = 0.2925937299273087
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 7.457778999999999
$ws.Range("N2").Value = 22.373337
$ws.Range("O2").Value = 0.08029647035915141
$ws.Range("P2").Value = 0.0802964703591514
$ws.Range("Q2").Value = 116.8574331927263
$ws.Range("R2").Value = 1051.716898734537
$ws.Range("S2").Value = 0.0234942437623817
$ws.Range("T2").Value = 0.02349424376238169

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 15.66920033333333
$ws.Range("H3").Value = 47.00760099999999
$ws.Range("I3").Value = 0.2925937299273087
$ws.Range("J3").Value = 0.2925937299273087
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 2.732509
$ws.Range("N3").Value = 8.197527
$ws.Range("O3").Value = 0.02942039820764526
$ws.Range("P3").Value = 0.02942039820764525
$ws.Range("Q3").Value = 42.81623093363633
$ws.Range("R3").Value = 385.3460784027269
$ws.Range("S3").Value = 0.008608224047521635
$ws.Range("T3").Value = 0.008608224047521632

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 15.66920033333333
$ws.Range("H4").Value = 47.00760099999999
$ws.Range("I4").Value = 0.2925937299273087
$ws.Range("J4").Value = 0.2925937299273087
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 82.68775466666666
$ws.Range("N4").Value = 248.063264
$ws.Range("O4").Value = 0.8902831314332034
$ws.Range("P4").Value = 0.8902831314332033
$ws.Range("Q4").Value = 1295.650992985518
$ws.Range("R4").Value = 11660.85893686966
$ws.Range("S4").Value = 0.2604912621174054
$ws.Range("T4").Value = 0.2604912621174054

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 16.47676966666667
$ws.Range("H5").Value = 49.43030900000001
$ws.Range("I5").Value = 0.3076736139282969
$ws.Range("J5").Value = 0.3076736139282968
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 7.457778999999999
$ws.Range("N5").Value = 22.373337
$ws.Range("O5").Value = 0.08029647035915141
$ws.Range("P5").Value = 0.0802964703591514
$ws.Range("Q5").Value = 122.8801068079037
$ws.Range("R5").Value = 1105.920961271133
$ws.Range("S5").Value = 0.02470510522108649
$ws.Range("T5").Value = 0.02470510522108648

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 16.47676966666667
$ws.Range("H6").Value = 49.43030900000001
$ws.Range("I6").Value = 0.3076736139282969
$ws.Range("J6").Value = 0.3076736139282968
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 2.732509
$ws.Range("N6").Value = 8.197527
$ws.Range("O6").Value = 0.02942039820764526
$ws.Range("P6").Value = 0.02942039820764525
$ws.Range("Q6").Value = 45.02292140509367
$ws.Range("R6").Value = 405.206292645843
$ws.Range("S6").Value = 0.009051880239755807
$ws.Range("T6").Value = 0.009051880239755802

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 16.47676966666667
$ws.Range("H7").Value = 49.43030900000001
$ws.Range("I7").Value = 0.3076736139282969
$ws.Range("J7").Value = 0.3076736139282968
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 82.68775466666666
$ws.Range("N7").Value = 248.063264
$ws.Range("O7").Value = 0.8902831314332034
$ws.Range("P7").Value = 0.8902831314332033
$ws.Range("Q7").Value = 1362.427087896509
$ws.Range("R7").Value = 12261.84379106858
$ws.Range("S7").Value = 0.2739166284674546
$ws.Range("T7").Value = 0.2739166284674545

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 21.406785
$ws.Range("H8").Value = 64.220355
$ws.Range("I8").Value = 0.3997326561443945
$ws.Range("J8").Value = 0.3997326561443944
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 7.457778999999999
$ws.Range("N8").Value = 22.373337
$ws.Range("O8").Value = 0.08029647035915141
$ws.Range("P8").Value = 0.0802964703591514
$ws.Range("Q8").Value = 159.647071630515
$ws.Range("R8").Value = 1436.823644674635
$ws.Range("S8").Value = 0.03209712137568323
$ws.Range("T8").Value = 0.03209712137568322

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 21.406785
$ws.Range("H9").Value = 64.220355
$ws.Range("I9").Value = 0.3997326561443945
$ws.Range("J9").Value = 0.3997326561443944
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 2.732509
$ws.Range("N9").Value = 8.197527
$ws.Range("O9").Value = 0.02942039820764526
$ws.Range("P9").Value = 0.02942039820764525
$ws.Range("Q9").Value = 58.49423267356499
$ws.Range("R9").Value = 526.4480940620849
$ws.Range("S9").Value = 0.01176029392036782
$ws.Range("T9").Value = 0.01176029392036782

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 21.406785
$ws.Range("H10").Value = 64.220355
$ws.Range("I10").Value = 0.3997326561443945
$ws.Range("J10").Value = 0.3997326561443944
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 82.68775466666666
$ws.Range("N10").Value = 248.063264
$ws.Range("O10").Value = 0.8902831314332034
$ws.Range("P10").Value = 0.8902831314332033
$ws.Range("Q10").Value = 1770.07898628208
$ws.Range("R10").Value = 15930.71087653872
$ws.Range("S10").Value = 0.3558752408483434
$ws.Range("T10").Value = 0.3558752408483433

